# Update countries & provincias Spain
#
# The underlying COVID dashboard source data refreshed: a handful of
# countries got new totals, and the table (sorted descending by "Casos
# totales", column B) was re-sorted to reflect the new ranking. A couple
# of countries tied on column B also swapped places relative to each
# other (no value change there - just re-fetch order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-CountryRow($name) {
    $found = $ws.Range("A4:A219").Find($name)
    if ($found -eq $null) {
        throw "Country not found: $name"
    }
    return $found.Row
}

# --- 1. Update the raw per-country figures that changed ------------------

$rArmenia = Get-CountryRow("Armenia")
$ws.Cells.Item($rArmenia, 2).Value = 7402
$ws.Cells.Item($rArmenia, 3).Value = 289
$ws.Cells.Item($rArmenia, 4).Value = 3220
$ws.Cells.Item($rArmenia, 5).Value = 4091
$ws.Cells.Item($rArmenia, 6).Value = 0
$ws.Cells.Item($rArmenia, 7).Value = 4
$ws.Cells.Item($rArmenia, 8).Value = 91

$rIsrael = Get-CountryRow("Israel")
$ws.Cells.Item($rIsrael, 2).Value = 16743
$ws.Cells.Item($rIsrael, 3).Value = 9
$ws.Cells.Item($rIsrael, 4).Value = 14362
$ws.Cells.Item($rIsrael, 5).Value = 2100

$rRumania = Get-CountryRow("Rumania")
$ws.Cells.Item($rRumania, 5).Value = 5446
$ws.Cells.Item($rRumania, 7).Value = 2
$ws.Cells.Item($rRumania, 8).Value = 1207

$rLetonia = Get-CountryRow("Letonia")
$ws.Cells.Item($rLetonia, 2).Value = 1053
$ws.Cells.Item($rLetonia, 3).Value = 4
$ws.Cells.Item($rLetonia, 4).Value = 741
$ws.Cells.Item($rLetonia, 5).Value = 290

# --- 2. Re-sort the whole table by "Casos totales" (column B) desc ------

$sortRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$sortRange.Sort($sortKey, 2)

# --- 3. Countries tied on column B that swapped positions on refresh ----

function Swap-Rows($nameA, $nameB) {
    $ra = Get-CountryRow($nameA)
    $rb = Get-CountryRow($nameB)
    $rngA = $ws.Range("A" + $ra + ":H" + $ra)
    $rngB = $ws.Range("A" + $rb + ":H" + $rb)
    $tmp = $rngA.Value2
    $rngA.Value2 = $rngB.Value2
    $rngB.Value2 = $tmp
}

Swap-Rows "Fiyi" "Curazao"
Swap-Rows "Santa Lucia" "Nueva Caledonia"
Swap-Rows "San Bartolome" "Bonaire, San Eustaquio y Saba"

# --- 4. Refresh the "last updated" timestamp caption ---------------------

$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 09:05"
